$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "29.483.55"
$ws.Range("E2").Value = "  +0.46%  "
$ws.Range("D3").Value = "1.878.76"
$ws.Range("E3").Value = "  +0.43%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.9996"
$ws.Range("D4").ClearFormats()
$ws.Range("E4").Value = "  -0.16%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "0.7150"
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = "  +0.36%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "242.28"
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = "  +0.47%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.9998"
$ws.Range("D7").ClearFormats()
$ws.Range("E7").Value = "  -0.18%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3119"
$ws.Range("D8").ClearFormats()
$ws.Range("E8").Value = "  +1.05%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.07733"
$ws.Range("D9").ClearFormats()
$ws.Range("E9").Value = "  -1.92%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "25.41"
$ws.Range("D10").ClearFormats()
$ws.Range("E10").Value = "  +0.08%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.08408"
$ws.Range("D11").ClearFormats()
$ws.Range("E11").Value = "  +1.86%  "
$ws.Range("D12").Value = "1.896.24"
$ws.Range("E12").Value = "  +1.49%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "5.268"
$ws.Range("D13").ClearFormats()
$ws.Range("E13").Value = "  +0.51%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.7197"
$ws.Range("D14").ClearFormats()
$ws.Range("E14").Value = "  -0.44%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "91.83"
$ws.Range("D15").ClearFormats()
$ws.Range("D16").Value = "29.506.82"
$ws.Range("E16").Value = "  +0.49%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.000008223"
$ws.Range("D17").ClearFormats()
$ws.Range("E17").Value = "  +5.14%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "6.001"
$ws.Range("D18").ClearFormats()
$ws.Range("E18").Value = "  +2.79%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "244.88"
$ws.Range("D19").ClearFormats()
$ws.Range("E19").Value = "  +0.33%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "13.28"
$ws.Range("D20").ClearFormats()
$ws.Range("E20").Value = "  +0.57%  "
$ws.Range("D21").Value = "2.132.50"
$ws.Range("E21").Value = "  -0.22%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.9995"
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = "  -0.19%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "7.954"
$ws.Range("D23").ClearFormats()
$ws.Range("E23").Value = "  -0.89%  "
$ws.Range("E24").Value = "  -0.18%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.1626"
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = "  +1.66%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "163.97"
$ws.Range("D26").ClearFormats()
$ws.Range("E26").Value = "  +0.97%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "9.049"
$ws.Range("D27").ClearFormats()
$ws.Range("E27").Value = "  +0.59%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "18.65"
$ws.Range("D28").ClearFormats()
$ws.Range("E28").Value = "  +2.26%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "1.513"
$ws.Range("D29").ClearFormats()
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "4.426"
$ws.Range("D30").ClearFormats()
$ws.Range("E30").Value = "  +1.07%  "
$ws.Range("E31").Value = "  -4.00%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "4.321"
$ws.Range("D32").ClearFormats()
$ws.Range("E32").Value = "  +5.26%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.05251"
$ws.Range("D33").ClearFormats()
$ws.Range("E33").Value = "  +1.26%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.936"
$ws.Range("D34").ClearFormats()
$ws.Range("E34").Value = "  -0.18%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.7773"
$ws.Range("D35").ClearFormats()
$ws.Range("E35").Value = "  +7.72%  "
$ws.Range("E36").Value = "  -0.71%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "2.677"
$ws.Range("D37").ClearFormats()
$ws.Range("E37").Value = "  +0.19%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.01869"
$ws.Range("D38").ClearFormats()
$ws.Range("E38").Value = "  +0.66%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "2.728"
$ws.Range("D39").ClearFormats()
$ws.Range("E39").Value = "  +1.24%  "
$ws.Range("D40").Value = "1.170.83"
$ws.Range("E40").Value = "  -0.36%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "6.439"
$ws.Range("D41").ClearFormats()
$ws.Range("E41").Value = "  +5.21%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "73.82"
$ws.Range("D42").ClearFormats()
$ws.Range("E42").Value = "  +1.64%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.8930"
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = "  -1.12%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "104.59"
$ws.Range("D44").ClearFormats()
$ws.Range("E44").Value = "  +2.41%  "
$ws.Range("E45").Value = "  -0.22%  "
$ws.Range("D46").Value = "2.031.57"
$ws.Range("E46").Value = "  +0.86%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "1.806"
$ws.Range("D47").ClearFormats()
$ws.Range("E47").Value = "  +0.95%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.5202"
$ws.Range("D48").ClearFormats()
$ws.Range("E48").Value = "  -1.66%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "9.429"
$ws.Range("D49").ClearFormats()
$ws.Range("E49").Value = "  +1.76%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.4325"
$ws.Range("D50").ClearFormats()
$ws.Range("E50").Value = "  +1.06%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "7.098"
$ws.Range("D51").ClearFormats()
$ws.Range("E51").Value = "  +1.12%  "
